$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H ("Status"), shifting the old "user"/"password"
# columns (and their column-width formatting) one place to the right.
$ws.Columns("H:H").Insert()

# Header row: rename the shifted columns and add the new Status header.
$ws.Range("H1").Value = "Status"
$ws.Range("I1").Value = "User"
$ws.Range("J1").Value = "Password"

# Row 2 (HTTP s1): target host corrected + service enabled.
$ws.Range("D2").Value = "shorter.rivetweb.org"
$ws.Range("H2").Value = "enable"

# Rows 3-5 (SSH / PGSQL for s1): disabled.
$ws.Range("H3").Value = "disable"
$ws.Range("H4").Value = "disable"
$ws.Range("H5").Value = "disable"

# Row 6 (HTTP s2): target host corrected + service enabled.
$ws.Range("D6").Value = "pastorious.rivetweb.org"
$ws.Range("H6").Value = "enable"

# Rows 7-8 (SSH / PGSQL for s2): disabled.
$ws.Range("H7").Value = "disable"
$ws.Range("H8").Value = "disable"

# View: zoom to 140% and move the selection to D4.
$excel.ActiveWindow.Zoom = 140
$ws.Range("D4").Select()
